$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the price / 1h-volume columns with freshly scraped values.
# For cells whose new text looks like a plain number (e.g. "1.00", "6.92")
# force the cell to Text format first so Excel keeps storing the exact
# original string instead of silently converting it to a number.
$ws.Range("D2").Value = '59.840.04'
$ws.Range("E2").Value = '  +0.10%  '
$ws.Range("D3").Value = '2.534.33'
$ws.Range("E3").Value = '  +1.37%  '
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '544.08'
$ws.Range("E5").Value = '  -0.08%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '145.71'
$ws.Range("E6").Value = '  -1.41%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.994'
$ws.Range("E7").Value = '  -0.11%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.573'
$ws.Range("E8").Value = '  -1.60%  '
$ws.Range("D9").Value = '2.562.43'
$ws.Range("E9").Value = '  +1.37%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.102'
$ws.Range("E10").Value = '  +0.31%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.161'
$ws.Range("E11").Value = '  +0.89%  '
$ws.Range("E12").Value = '  +2.84%  '
$ws.Range("E13").Value = '  +0.97%  '
$ws.Range("D14").Value = '2.978.97'
$ws.Range("E14").Value = '  +1.07%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '23.59'
$ws.Range("E15").Value = '  -4.51%  '
$ws.Range("D16").Value = '59.770.23'
$ws.Range("E16").Value = '  +0.00%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000144'
$ws.Range("E17").Value = '  +2.31%  '
$ws.Range("D18").Value = '2.551.64'
$ws.Range("E18").Value = '  +1.37%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.25'
$ws.Range("E19").Value = '  -3.03%  '
$ws.Range("E20").Value = '  -1.61%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '327.94'
$ws.Range("E21").Value = '  +0.06%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.00'
$ws.Range("E22").Value = '  +0.09%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.95'
$ws.Range("E23").Value = '  +1.83%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '62.16'
$ws.Range("E24").Value = '  +1.13%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.442'
$ws.Range("E25").Value = '  -1.51%  '
$ws.Range("E26").Value = '  +1.76%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.992'
$ws.Range("E27").Value = '  -1.80%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.02'
$ws.Range("E28").Value = '  +0.95%  '
$ws.Range("B29").Value = 'PEPE'
$ws.Range("C29").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D29").Value = '0.0₃0800'
$ws.Range("E29").Value = '  +0.31%  '
$ws.Range("B30").Value = 'Aptos'
$ws.Range("C30").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.92'
$ws.Range("E30").Value = '  -2.88%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.82'
$ws.Range("E31").Value = '  -1.01%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.21'
$ws.Range("E32").Value = '  -7.04%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '162.54'
$ws.Range("E33").Value = '  +2.37%  '
$ws.Range("E34").Value = '  +2.27%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.998'
$ws.Range("E35").Value = '  -0.08%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '18.83'
$ws.Range("E36").Value = '  -0.53%  '
$ws.Range("E37").Value = '  -1.25%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.64'
$ws.Range("E38").Value = '  -5.89%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.74'
$ws.Range("E39").Value = '  -5.92%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '37.19'
$ws.Range("E40").Value = '  +1.14%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '303.29'
$ws.Range("E41").Value = '  -4.60%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.839'
$ws.Range("E42").Value = '  +0.02%  '
$ws.Range("E43").Value = '  -1.98%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.993'
$ws.Range("E44").Value = '  -0.09%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.607'
$ws.Range("E45").Value = '  +0.34%  '
$ws.Range("E46").Value = '  +0.27%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '19.03'
$ws.Range("E47").Value = '  +1.54%  '
$ws.Range("E48").Value = '  -0.85%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '124.15'
$ws.Range("E49").Value = '  -2.58%  '
$ws.Range("E50").Value = '  -2.57%  '
$ws.Range("E51").Value = '  -1.60%  '
